$d = $word.ActiveDocument
$replacements = @(
    @("2024-07-02 Tuesday", "2024-07-03 Wednesday"),
    @("281÷8=35, 1", "965÷3=321, 2"),
    @("689÷8=86, 1", "647÷9=71, 8"),
    @("676÷2=338, 0", "964÷5=192, 4"),
    @("746÷6=124, 2", "931÷8=116, 3"),
    @("998÷2=499, 0", "671÷6=111, 5"),
    @("893÷8=111, 5", "343÷3=114, 1"),
    @("150÷9=16, 6", "380÷9=42, 2"),
    @("863÷6=143, 5", "208÷8=26, 0"),
    @("992÷7=141, 5", "141÷9=15, 6"),
    @("973÷4=243, 1", "883÷6=147, 1"),
    @("731÷7=104, 3", "457÷2=228, 1"),
    @("197÷2=98, 1", "909÷7=129, 6"),
    @("311÷3=103, 2", "492÷9=54, 6"),
    @("428÷9=47, 5", "147÷2=73, 1"),
    @("884÷2=442, 0", "446÷7=63, 5"),
    @("523÷4=130, 3", "208÷6=34, 4"),
    @("843÷6=140, 3", "771÷4=192, 3"),
    @("516÷8=64, 4", "547÷7=78, 1"),
    @("319÷4=79, 3", "400÷3=133, 1"),
    @("247÷8=30, 7", "795÷2=397, 1"),
    @("134÷2=67, 0", "989÷7=141, 2"),
    @("685÷5=137, 0", "594÷9=66, 0"),
    @("220÷5=44, 0", "494÷5=98, 4"),
    @("600÷6=100, 0", "707÷7=101, 0"),
    @("615÷3=205, 0", "272÷4=68, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: text not found: $old"
    }
}

Write-Host "Replacements complete."
